$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($range, $value) {
    # Force the cell to keep a text data type (matching the source data,
    # which stores prices as text) rather than letting Excel auto-convert
    # numeric-looking strings into numbers.
    $range.NumberFormat = "@"
    $range.Value = $value
}

# Row 2 - Bitcoin
Set-TextValue $ws.Range("D2") "30.414.90"
$ws.Range("E2").Value = "  +2.58%  "

# Row 3 - Ethereum
Set-TextValue $ws.Range("D3") "2.110.17"
$ws.Range("E3").Value = "  +0.62%  "

# Row 4 - TetherUSD
Set-TextValue $ws.Range("D4") "1.006"
$ws.Range("E4").Value = "  -0.12%  "

# Row 5 - BNB
Set-TextValue $ws.Range("D5") "345.14"
$ws.Range("E5").Value = "  +0.77%  "

# Row 6 - USDC
$ws.Range("E6").Value = "  -0.02%  "

# Row 7 - XRP
Set-TextValue $ws.Range("D7") "0.5233"
$ws.Range("E7").Value = "  +2.13%  "

# Row 8 - Cardano
Set-TextValue $ws.Range("D8") "0.4448"
$ws.Range("E8").Value = "  +1.10%  "

# Row 9 - OKB
Set-TextValue $ws.Range("D9") "54.53"
$ws.Range("E9").Value = "  +2.47%  "

# Row 10 - Dogecoin
Set-TextValue $ws.Range("D10") "0.09392"
$ws.Range("E10").Value = "  +2.76%  "

# Row 11 - Polygon
Set-TextValue $ws.Range("D11") "1.175"
$ws.Range("E11").Value = "  +0.36%  "

# Row 12 - Solana
Set-TextValue $ws.Range("D12") "25.00"
$ws.Range("E12").Value = "  +0.83%  "

# Row 13 - Chainlink
Set-TextValue $ws.Range("D13") "8.682"
$ws.Range("E13").Value = "  +6.06%  "

# Row 14 - Polkadot
Set-TextValue $ws.Range("D14") "6.960"
$ws.Range("E14").Value = "  +3.20%  "

# Row 15 - WrappedEther
Set-TextValue $ws.Range("D15") "2.113.12"
$ws.Range("E15").Value = "  +0.74%  "

# Row 16 - Litecoin
Set-TextValue $ws.Range("D16") "101.96"
$ws.Range("E16").Value = "  +2.24%  "

# Row 17 - ShibaInu
Set-TextValue $ws.Range("D17") "0.00001164"
$ws.Range("E17").Value = "  +1.68%  "

# Row 18 - BinanceUSD
Set-TextValue $ws.Range("D18") "1.007"
$ws.Range("E18").Value = "  +0.00%  "

# Row 19 - Avalanche
Set-TextValue $ws.Range("D19") "21.24"
$ws.Range("E19").Value = "  +0.73%  "

# Row 20 - TRON
Set-TextValue $ws.Range("D20") "0.06730"
$ws.Range("E20").Value = "  +1.35%  "

# Row 21 - Uniswap
Set-TextValue $ws.Range("D21") "6.353"
$ws.Range("E21").Value = "  +2.81%  "

# Row 22 - Dai
$ws.Range("E22").Value = "  -0.12%  "

# Row 23 - WrappedBTC
Set-TextValue $ws.Range("D23") "30.415.02"
$ws.Range("E23").Value = "  +2.41%  "

# Row 24 - Cosmos
Set-TextValue $ws.Range("D24") "12.65"
$ws.Range("E24").Value = "  +0.58%  "

# Row 25 - Toncoin
Set-TextValue $ws.Range("D25") "2.302"
$ws.Range("E25").Value = "  -0.20%  "

# Row 26 - EthereumClassic
Set-TextValue $ws.Range("D26") "22.06"
$ws.Range("E26").Value = "  +1.10%  "

# Row 27 - was Monero, now LidoDAOToken
$ws.Range("B27").Value = "LidoDAOToken"
$ws.Range("C27").Value = "https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo"
Set-TextValue $ws.Range("D27") "2.541"
$ws.Range("E27").Value = "  +0.70%  "

# Row 28 - was LidoDAOToken, now Monero
$ws.Range("B28").Value = "Monero"
$ws.Range("C28").Value = "https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr"
Set-TextValue $ws.Range("D28") "163.23"
$ws.Range("E28").Value = "  +0.58%  "

# Row 29 - BitcoinCash
Set-TextValue $ws.Range("D29") "134.07"
$ws.Range("E29").Value = "  +1.23%  "

# Row 30 - ImmutableX
Set-TextValue $ws.Range("D30") "1.157"
$ws.Range("E30").Value = "  +2.36%  "

# Row 31 - ARBITRUM
Set-TextValue $ws.Range("D31") "1.740"
$ws.Range("E31").Value = "  +6.35%  "

# Row 32 - Stellar
Set-TextValue $ws.Range("D32") "0.1055"
$ws.Range("E32").Value = "  +1.08%  "

# Row 33 - InternetComputer(DFINITY)
Set-TextValue $ws.Range("D33") "6.829"
$ws.Range("E33").Value = "  +13.03%  "

# Row 34 - Filecoin
Set-TextValue $ws.Range("D34") "6.279"
$ws.Range("E34").Value = "  +2.01%  "

# Row 35 - HuobiToken
$ws.Range("E35").Value = "  -0.97%  "

# Row 36 - FraxShare
$ws.Range("E36").Value = "  +1.58%  "

# Row 37 - VeChain
Set-TextValue $ws.Range("D37") "0.02627"
$ws.Range("E37").Value = "  +2.34%  "

# Row 38 - Hedera
Set-TextValue $ws.Range("D38") "0.06799"
$ws.Range("E38").Value = "  +2.09%  "

# Row 39 - TheSandbox
Set-TextValue $ws.Range("D39") "0.7064"
$ws.Range("E39").Value = "  +3.19%  "

# Row 40 - TrustWalletToken
Set-TextValue $ws.Range("D40") "1.360"
$ws.Range("E40").Value = "  +5.70%  "

# Row 41 - Aptos
Set-TextValue $ws.Range("D41") "12.59"
$ws.Range("E41").Value = "  +1.84%  "

# Row 42 - Algorand
Set-TextValue $ws.Range("D42") "0.2228"
$ws.Range("E42").Value = "  -0.28%  "

# Row 43 - Decentraland
Set-TextValue $ws.Range("D43") "0.6863"
$ws.Range("E43").Value = "  +2.94%  "

# Row 44 - EnergySwap
Set-TextValue $ws.Range("D44") "14.41"
$ws.Range("E44").Value = "  +1.29%  "

# Row 45 - NEARProtocol
Set-TextValue $ws.Range("D45") "2.362"

# Row 46 - Frax
Set-TextValue $ws.Range("D46") "1.006"
$ws.Range("E46").Value = "  +0.03%  "

# Row 47 - WEMIXTOKEN
Set-TextValue $ws.Range("D47") "1.365"
$ws.Range("E47").Value = "  +17.48%  "

# Row 48 - PancakeSwap
Set-TextValue $ws.Range("D48") "3.648"
$ws.Range("E48").Value = "  +1.10%  "

# Row 49 - BabyDogeCoin
Set-TextValue $ws.Range("D49") "0.00000000347"
$ws.Range("E49").Value = "  +4.14%  "

# Row 50 - ThetaToken
Set-TextValue $ws.Range("D50") "1.209"
$ws.Range("E50").Value = "  +9.36%  "

# Row 51 - EOS
Set-TextValue $ws.Range("D51") "1.221"
$ws.Range("E51").Value = "  +0.24%  "
